# Weekly update: two new "Ciboulette" price records were collected and
# inserted into the historical series, pushing the existing rows (which are
# ordered most-recent-first) down to make room. One record is inserted right
# before the old row 209, and a second is inserted a bit further down the
# series (right before what was row 313 originally, i.e. row 314 after the
# first insert). Everything else in the sheet is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common / constant values shared by every data row in this sheet.
$mercadoId   = 3
$mercado     = "Femacal de La Calera"
$region      = "Coquimbo"
$codreg      = 5
$categoriaId = 100112039
$categoria   = "Ciboulette"
$variedad    = "Sin especificar"
$calidad     = "Primera"
$unidad      = "`$/docena de atados"
$origen      = "Provincia de Quillota"
$kgUnidades  = 3
$clasif      = "Hortaliza"

function Set-DataRow {
    param($row, $fecha, $volumen, $precioMin, $precioMax, $precioProm, $precioKg)

    $ws.Range("A$row").Value = $mercadoId
    $ws.Range("B$row").Value = $mercado
    $ws.Range("C$row").Value = $region
    $ws.Range("D$row").Value = $fecha
    $ws.Range("E$row").Value = $codreg
    $ws.Range("F$row").Value = $categoriaId
    $ws.Range("G$row").Value = $categoria
    $ws.Range("H$row").Value = $variedad
    $ws.Range("I$row").Value = $calidad
    $ws.Range("J$row").Value = $volumen
    $ws.Range("K$row").Value = $precioMin
    $ws.Range("L$row").Value = $precioMax
    $ws.Range("M$row").Value = $precioProm
    $ws.Range("N$row").Value = $unidad
    $ws.Range("O$row").Value = $origen
    $ws.Range("P$row").Value = $precioKg
    $ws.Range("Q$row").Value = $kgUnidades
    $ws.Range("R$row").Value = $clasif
}

# --- First insertion: brand-new row at (old/new) position 209 -------------
# Shifts old rows 209..312 down to 210..313.
$ws.Rows.Item(209).Insert()
Set-DataRow -row 209 -fecha 44818 -volumen 160 -precioMin 1500 -precioMax 1500 -precioProm 1500 -precioKg 500

# --- Second insertion: brand-new row at position 314 -----------------------
# (This is old row 313's new position before this second insert; inserting
# here shifts old rows 313..363 -- now sitting at 314..364 -- down again to
# 315..365.)
$ws.Rows.Item(314).Insert()
Set-DataRow -row 314 -fecha 44816 -volumen 120 -precioMin 1500 -precioMax 1500 -precioProm 1500 -precioKg 500
